# RF classify (read data header bug)
#
# The sheet used to keep a duplicate "raw" copy of the header/timestamp
# columns in H:M (computed from a video-start-time offset in B9/B11) and
# then re-derived the real A:B start/end times from those raw columns via
# formulas (e.g. A2 = H2-B$9). That raw header block was a bug/leftover
# from the read-data step, so it gets removed: the A:B columns are frozen
# to their already-computed values, the H:M helper block is cleared, and
# the now-unused "Video start time:" label row is deleted (shifting the
# trailing helper row up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Freeze A2:B7 (currently formulas referencing the raw H:M block and B9)
# down to plain values - the formulas are the "header bug" being removed.
$ws.Range("A2:B7").Value = $ws.Range("A2:B7").Value2

# Drop the duplicated raw header/data block in columns H:M (rows 1-7).
$ws.Range("H1:M7").Clear()

# Remove the now unused "Video start time:" row; this shifts the former
# row 11 (the B11 helper constant) up to row 10.
$ws.Rows("9:9").Delete()

# Match the author's final selection state.
[void]$ws.Range("K11").Select()
